$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 so that the existing data (rows 2-12) shifts
# down to rows 3-13, then fill in the new row 2 with the new weekly record.
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits the formatting of the row above it
# (the bold header row). Reset it back to the plain default style used by
# the rest of the data rows, then re-apply the date number format on the
# "Fecha" column only, matching the rest of the table.
$ws.Rows.Item(2).ClearFormats()
$ws.Cells.Item(3, 4).Copy()
$ws.Cells.Item(2, 4).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(2, 3).Value = "Bíobío"
$ws.Cells.Item(2, 4).Value = 44860
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = 300000000
$ws.Cells.Item(2, 7).Value = "Espárragos"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 1100
$ws.Cells.Item(2, 11).Value = 1500
$ws.Cells.Item(2, 12).Value = 1700
$ws.Cells.Item(2, 13).Value = 1609
$ws.Cells.Item(2, 14).Value = "`$/kilo"
$ws.Cells.Item(2, 15).Value = "Provincia de Linares"
$ws.Cells.Item(2, 16).Value = 1609
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = "Hortaliza"
